# Update the Supplier for the "Blackfly BFS-U3-51S5" camera row
# from "PointGrey" to "FLIR" (PointGrey was acquired by and rebranded as FLIR).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C3").Value = "FLIR"
